# Daily attendance processing - 2025-11-07 09:46:07
# Updates the "Recorded By" (column G) values on the active worksheet,
# reordering the contributor list so that "System"/"system" is listed last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [PSCustomObject]@{ Row = 2; Value = "backup@backdoor.com, system, System" }
    [PSCustomObject]@{ Row = 3; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 4; Value = "backup@backdoor.com, System" }
    [PSCustomObject]@{ Row = 6; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 7; Value = "admin@admin.com, System" }
    [PSCustomObject]@{ Row = 10; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 12; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 13; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 14; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 15; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 18; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 19; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 20; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 21; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 22; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 24; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 26; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 28; Value = "backup@backdoor.com, system, System" }
    [PSCustomObject]@{ Row = 29; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 30; Value = "backup@backdoor.com, System" }
    [PSCustomObject]@{ Row = 32; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 33; Value = "admin@admin.com, System" }
    [PSCustomObject]@{ Row = 36; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 38; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 39; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 40; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 41; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 44; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 45; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 46; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 47; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 48; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 50; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 52; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 54; Value = "backup@backdoor.com, system, System" }
    [PSCustomObject]@{ Row = 55; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 56; Value = "backup@backdoor.com, System" }
    [PSCustomObject]@{ Row = 58; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 59; Value = "admin@admin.com, System" }
    [PSCustomObject]@{ Row = 62; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 64; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 65; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 66; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 67; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 70; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 71; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 72; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 73; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 74; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 76; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 78; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 83; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 84; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 85; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 86; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 90; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 92; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 99; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 101; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 109; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 110; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 111; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 112; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 116; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 118; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 125; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 127; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 135; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 136; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 137; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 138; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 142; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 144; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 151; Value = "dnasr281@gmail.com, System" }
    [PSCustomObject]@{ Row = 153; Value = "dnasr281@gmail.com, System" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}

Write-Output "Updated $($updates.Count) cells in column G (Recorded By)."
